$wb = $excel.ActiveWorkbook

# Update the DEA output/input labels on LP2..LP5 so each sales-office sheet
# gets its own x1n/y1n/y2n labels instead of the copy-pasted x11/y11/y21.
# Shared-string append order matches the captured edit: J4/K4 across all
# four sheets first, then I4 across all four sheets.

$ws2 = $wb.Worksheets.Item("LP2")
$ws2.Range("J4").Value = "y12"
$ws2.Range("K4").Value = "y22"

$ws3 = $wb.Worksheets.Item("LP3")
$ws3.Range("J4").Value = "y13"
$ws3.Range("K4").Value = "y23"

$ws4 = $wb.Worksheets.Item("LP4")
$ws4.Range("J4").Value = "y14"
$ws4.Range("K4").Value = "y24"

$ws5 = $wb.Worksheets.Item("LP5")
$ws5.Range("J4").Value = "y15"
$ws5.Range("K4").Value = "y25"

$ws2.Range("I4").Value = "x12"
$ws3.Range("I4").Value = "x13"
$ws4.Range("I4").Value = "x14"
$ws5.Range("I4").Value = "x15"

# Update the selection on each non-active sheet to I4 (matches the saved
# workbook state), then finish on LP5 so it becomes the active tab/sheet.
$ws2.Activate()
[void]$ws2.Range("I4").Select()

$ws3.Activate()
[void]$ws3.Range("I4").Select()

$ws4.Activate()
[void]$ws4.Range("I4").Select()

$ws5.Activate()
[void]$ws5.Range("I4").Select()
